$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark m_functionRemove_Type_Event (A18) as done by copying the
# strikethrough style already used for the other finished items (A16/A17).
$ws.Range("A16").Copy() | Out-Null
$ws.Range("A18").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# Remove the stale "see m_functionAdd_Type_Event" note next to it.
$ws.Range("E18").ClearContents()

# The TODO list header moves from "finish" to "test".
$ws.Range("A15").Value = "methods in code.js to test"

# "write code methods listed below" is also now done.
$ws.Range("A5").Value = "done"

# Restore the original selection.
$ws.Range("A16").Select() | Out-Null
